$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace en-dash with hyphen / reassign team names for column B
$ws.Range("B4").Value = 'Team Tine - Danfoss'
$ws.Range("B11").Value = 'Konami - Yokohama'
$ws.Range("B15").Value = 'Fanox - Bizkaia'
$ws.Range("B19").Value = 'Ciudad de Medellin'
$ws.Range("B20").Value = 'Saxo Bank - Sungard'
$ws.Range("B21").Value = 'Manuela Fundacion'
$ws.Range("B22").Value = 'Adriatic.hr'
$ws.Range("B23").Value = 'Barclays - Engie'
$ws.Range("B24").Value = 'Nutella Pro Team'
$ws.Range("B25").Value = 'Boels - Dolmans'
$ws.Range("B26").Value = 'Acqua & Sapone'
$ws.Range("B27").Value = 'Omega Pharma'
$ws.Range("B28").Value = 'Frigo - Philips'
$ws.Range("B29").Value = 'OTIS Cycling Team'
$ws.Range("B30").Value = 'Saunier Duval'
$ws.Range("B31").Value = 'Tissot'
$ws.Range("B32").Value = 'Televisión'

# Renumber ID column A for rows 19-32
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29
$ws.Range("A31").Value = 30
$ws.Range("A32").Value = 31

# Update selection
$ws.Range("B29").Select()
